$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.824.37"
$ws.Range("E2").Value = "  +1.35%  "
$ws.Range("D3").Value = "2.426.53"
$ws.Range("E3").Value = "  +1.66%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'307.29"
$ws.Range("E5").Value = "  +1.20%  "
$ws.Range("D6").Value = "'97.36"
$ws.Range("E6").Value = "  -0.38%  "
$ws.Range("D7").Value = "'0.511"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "'0.494"
$ws.Range("E9").Value = "  -1.87%  "
$ws.Range("D10").Value = "'35.31"
$ws.Range("E10").Value = "  +2.63%  "
$ws.Range("D11").Value = "'0.0798"
$ws.Range("E11").Value = "  +0.86%  "
$ws.Range("E12").Value = "  +2.23%  "
$ws.Range("D13").Value = "'18.54"
$ws.Range("E13").Value = "  -0.17%  "
$ws.Range("D14").Value = "'6.90"
$ws.Range("E14").Value = "  +1.37%  "
$ws.Range("D15").Value = "2.797.72"
$ws.Range("E15").Value = "  +1.72%  "
$ws.Range("D16").Value = "2.418.85"
$ws.Range("E16").Value = "  +1.69%  "
$ws.Range("D17").Value = "'0.830"
$ws.Range("E17").Value = "  +2.28%  "
$ws.Range("D18").Value = "43.836.82"
$ws.Range("E18").Value = "  +1.37%  "
$ws.Range("E19").Value = "  +0.88%  "
$ws.Range("D20").Value = "'12.13"
$ws.Range("E20").Value = "  -1.75%  "
$ws.Range("D21").Value = "0.0₃0902"
$ws.Range("E21").Value = "  +1.35%  "
$ws.Range("D22").Value = "'68.16"
$ws.Range("E22").Value = "  -0.34%  "
$ws.Range("D23").Value = "'238.55"
$ws.Range("E23").Value = "  +0.85%  "
$ws.Range("D24").Value = "'2.27"
$ws.Range("E24").Value = "  +1.74%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").Value = "'2.46"
$ws.Range("E26").Value = "  +0.53%  "
$ws.Range("D27").Value = "'25.06"
$ws.Range("E27").Value = "  +0.55%  "
$ws.Range("E28").Value = "  -0.56%  "
$ws.Range("D29").Value = "'9.46"
$ws.Range("E29").Value = "  +3.20%  "
$ws.Range("D30").Value = "'32.40"
$ws.Range("E30").Value = "  +2.12%  "
$ws.Range("E31").Value = "  +17.98%  "
$ws.Range("D32").Value = "'18.52"
$ws.Range("E32").Value = "  +7.02%  "
$ws.Range("D33").Value = "'5.15"
$ws.Range("E33").Value = "  +0.27%  "
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("D35").Value = "'0.0756"
$ws.Range("E35").Value = "  +3.69%  "
$ws.Range("D36").Value = "'1.92"
$ws.Range("E36").Value = "  +2.86%  "
$ws.Range("D37").Value = "'130.53"
$ws.Range("E37").Value = "  +25.20%  "
$ws.Range("E38").Value = "  +3.83%  "
$ws.Range("D39").Value = "'4.42"
$ws.Range("E39").Value = "  +0.45%  "
$ws.Range("E40").Value = "  -1.60%  "
$ws.Range("D41").Value = "'0.109"
$ws.Range("E41").Value = "  -0.20%  "
$ws.Range("D42").Value = "'21.26"
$ws.Range("E42").Value = "  -6.74%  "
$ws.Range("D43").Value = "1.947.36"
$ws.Range("E43").Value = "  -0.14%  "
$ws.Range("D44").Value = "'0.0284"
$ws.Range("E44").Value = "  +1.00%  "
$ws.Range("E45").Value = "  +2.15%  "
$ws.Range("D46").Value = "'2.85"
$ws.Range("E46").Value = "  +3.01%  "
$ws.Range("D47").Value = "'9.32"
$ws.Range("E47").Value = "  -0.64%  "
$ws.Range("D48").Value = "2.656.75"
$ws.Range("E48").Value = "  +1.99%  "
$ws.Range("D49").Value = "'1.58"
$ws.Range("E49").Value = "  +4.49%  "
$ws.Range("D50").Value = "'52.77"
$ws.Range("E50").Value = "  -0.77%  "
$ws.Range("D51").Value = "'72.69"
$ws.Range("E51").Value = "  +0.33%  "
